$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Требования к форме").Name = "Form requirements"
$wb.Worksheets.Item("Все поля").Name = "All fields"
